# Generate Report for Handback
# The handback report was regenerated: the "3a84e157..." file's handback is
# now detected as out of sync with en-US (status text flips for every row
# that shares that status string), and the "1e0e29bd..." file just got a
# fresh handback timestamp recorded on both the zh-cn and de-de report
# sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Flip the "in sync" status to "not in sync" everywhere it is shown
# (Overview zh-cn/de-de status columns, and the Status column on each
# language report sheet).
$overview.Range("E2:F3").Value = "Handed back: not in sync with en-US"
$zhcn.Range("C2:C3").Value = "Handed back: not in sync with en-US"
$dede.Range("C2:C3").Value = "Handed back: not in sync with en-US"

# Record the freshly generated handback datetimes for the 1e0e29bd row.
$zhcn.Range("K2").Value = "2016-10-20 01:02:35"
$dede.Range("K2").Value = "2016-10-20 01:02:53"

# The longer status text widened the status columns (matches the report
# generator's column auto-sizing for the new "not in sync" wording).
$overview.Columns.Item(5).ColumnWidth = 32.6
$overview.Columns.Item(6).ColumnWidth = 32.6
$zhcn.Columns.Item(3).ColumnWidth = 32.6
$dede.Columns.Item(3).ColumnWidth = 32.6
